$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "BVTs": restyle row 2 (drop the redundant per-font xf's that Excel
# collapsed back to the plain/default + wrap-only styles), then append the
# new "Drillthrough" BVT as row 20.
# ---------------------------------------------------------------------------
$bvts = $wb.Worksheets.Item("BVTs")

# Row 2 lost its (now-unused) applyFont-only styles; A2/B2/C2/E2 fall back to
# the default (no) style, D2 keeps the wrap-text-only style already used by
# sibling cells like D4/E3.
$bvts.Range("A3").Copy()
$bvts.Range("A2:C2").PasteSpecial(-4122)
$bvts.Range("E2").PasteSpecial(-4122)

$bvts.Range("E3").Copy()
$bvts.Range("D2").PasteSpecial(-4122)

# New row 20 - "Drillthrough" BVT (API 2.5 / context-menu support).
$bvts.Range("A20").Value = 19
$bvts.Range("B20").Value = "Drillthrough"

$bvts.Range("E19").Copy()
$bvts.Range("C20").PasteSpecial(-4122)
$bvts.Range("D20").PasteSpecial(-4122)
$bvts.Range("E20").PasteSpecial(-4122)

$bvts.Range("C20").Value = "Created custom menu to drill through from one visual to another."
$bvts.Range("D20").Value = "1.Generate a chart with some data" + [char]10 + "2.Create new report page and in DrillThrough add the fields for drillthrough." + [char]10 + "3. Right click on the chart, select the Drillthrough option from the menu. "
$bvts.Range("E20").Value = "1. On right click of the chart and selecting the drillthrough option from the context-menu , the report will drillthrough to the newly created report page."

$bvts.Rows.Item(20).RowHeight = 75

# ---------------------------------------------------------------------------
# View state: BVTs becomes the active/visible tab (was Checklist), with a new
# selection further down the new row; Checklist keeps its own selection but
# is no longer the tab in front.
# ---------------------------------------------------------------------------
$checklist = $wb.Worksheets.Item("Checklist")
$checklist.Activate()
$checklist.Range("G25").Select()

$bvts.Activate()
$bvts.Range("E22").Select()
